$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact_Information")

# Fix up the Enquiry text (column C) wording/punctuation
$ws.Range("C2").Value = "Help, I've fallen and I cannot get up!"
$ws.Range("C5").Value = "Please help with my order."
$ws.Range("C6").Value = "Great Products!"
$ws.Range("C7").Value = "Love the prices!"
$ws.Range("C8").Value = "Is shipping free?"
$ws.Range("C9").Value = "Seriously?"
$ws.Range("C10").Value = "I bought it all!!!"

# Update selection to match the new active cell
$ws.Range("C10").Select()
